$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Sun Jul 21 03:59:40 UTC 2024 with GitHub Actions
# Price (col D) and Volume(1h) (col E) refresh per row.
# A leading "'" forces numeric-looking price strings to stay text (matches
# the source data, which stores prices as plain strings, not numbers).
$ws.Range("D2").Value = "67.173.60"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.512.69"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'595.84"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "'173.16"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("E9").Value = "  +5.39%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "4.121.97"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'29.27"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").Value = "67.128.39"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "3.487.97"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "'14.17"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "'396.84"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'8.05"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "'73.08"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'10.26"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'6.33"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'1.46"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "'2.07"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "'7.39"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("D35").Value = "'163.90"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "'0.887"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'0.0753"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'26.60"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'27.21"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("D43").Value = "2.839.06"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("D45").Value = "'42.82"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "'0.0304"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").Value = "'339.52"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").Value = "'34.68"
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'6.51"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'0.850"
$ws.Range("E51").Value = "  -0.68%  "
